$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.538.85"
$ws.Range("E2").Value = "  +0.47%  "
Set-TextValue $ws.Range("D3") "1.794.78"
$ws.Range("E3").Value = "  -0.43%  "
Set-TextValue $ws.Range("D4") "0.9998"
$ws.Range("E4").Value = "  -1.05%  "
Set-TextValue $ws.Range("D5") "339.25"
$ws.Range("E5").Value = "  +1.67%  "
Set-TextValue $ws.Range("D6") "0.9965"
$ws.Range("E6").Value = "  -0.89%  "
Set-TextValue $ws.Range("D7") "0.3924"
$ws.Range("E7").Value = "  +3.53%  "
Set-TextValue $ws.Range("D8") "0.3468"
$ws.Range("E8").Value = "  -0.67%  "
Set-TextValue $ws.Range("D9") "48.33"
$ws.Range("E9").Value = "  -0.82%  "
Set-TextValue $ws.Range("D10") "1.198"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -1.10%  "
Set-TextValue $ws.Range("D12") "0.9975"
$ws.Range("E12").Value = "  -1.14%  "
Set-TextValue $ws.Range("D13") "21.98"
$ws.Range("E13").Value = "  +0.32%  "
Set-TextValue $ws.Range("D14") "6.523"
$ws.Range("E14").Value = "  +0.14%  "
Set-TextValue $ws.Range("D15") "1.793.91"
$ws.Range("E15").Value = "  -0.80%  "
Set-TextValue $ws.Range("D16") "7.171"
$ws.Range("E16").Value = "  +1.18%  "
Set-TextValue $ws.Range("D17") "0.00001101"
$ws.Range("E17").Value = "  -0.28%  "
Set-TextValue $ws.Range("D18") "0.06684"
$ws.Range("E18").Value = "  -0.13%  "
Set-TextValue $ws.Range("D19") "85.00"
$ws.Range("E19").Value = "  -0.23%  "
Set-TextValue $ws.Range("D20") "0.9969"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  +2.16%  "
Set-TextValue $ws.Range("D22") "6.571"
$ws.Range("E22").Value = "  +1.61%  "
Set-TextValue $ws.Range("D23") "27.536.50"
$ws.Range("E23").Value = "  +0.33%  "
Set-TextValue $ws.Range("D24") "12.47"
$ws.Range("E24").Value = "  -1.00%  "
Set-TextValue $ws.Range("D25") "2.407"
$ws.Range("E25").Value = "  -1.67%  "
Set-TextValue $ws.Range("D26") "2.515"
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D27") "1.470"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "21.25"
$ws.Range("E28").Value = "  -1.82%  "
Set-TextValue $ws.Range("D29") "156.47"
$ws.Range("E29").Value = "  +4.26%  "
Set-TextValue $ws.Range("D30") "1.998.31"
$ws.Range("E30").Value = "  -0.82%  "
Set-TextValue $ws.Range("D31") "135.38"
$ws.Range("E31").Value = "  +0.85%  "
Set-TextValue $ws.Range("D32") "4.043"
$ws.Range("E32").Value = "  -0.94%  "
Set-TextValue $ws.Range("D33") "6.066"
$ws.Range("E33").Value = "  -0.70%  "
Set-TextValue $ws.Range("D34") "0.08790"
$ws.Range("E34").Value = "  +1.19%  "
Set-TextValue $ws.Range("D35") "13.10"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D36") "5.475"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D37") "1.620"
$ws.Range("E37").Value = "  -4.02%  "
Set-TextValue $ws.Range("D38") "0.02417"
$ws.Range("E38").Value = "  +2.32%  "
Set-TextValue $ws.Range("D39") "0.06487"
$ws.Range("E39").Value = "  +1.52%  "
Set-TextValue $ws.Range("D40") "0.6842"
$ws.Range("E40").Value = "  -0.29%  "
Set-TextValue $ws.Range("D41") "0.2217"
$ws.Range("E41").Value = "  -0.06%  "
Set-TextValue $ws.Range("D42") "1.256"
$ws.Range("E42").Value = "  -2.57%  "
Set-TextValue $ws.Range("D43") "8.414"
$ws.Range("E43").Value = "  -6.23%  "
Set-TextValue $ws.Range("D44") "14.43"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.6410"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D46") "0.9965"
$ws.Range("E46").Value = "  -0.88%  "
Set-TextValue $ws.Range("D47") "3.876"
$ws.Range("E47").Value = "  +0.82%  "
Set-TextValue $ws.Range("D48") "2.140"
$ws.Range("E48").Value = "  +0.31%  "
Set-TextValue $ws.Range("D49") "132.30"
$ws.Range("E49").Value = "  +0.79%  "
Set-TextValue $ws.Range("D50") "0.07203"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  -0.03%  "
